$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the width of column F (to match for the two newly inserted columns).
$refWidth = $ws.Columns.Item(6).ColumnWidth

# Insert two new columns at H:I (everything from the old column H onward shifts right by 2).
$ws.Range("H:I").Insert()

# Give the new H,I columns the same width as the neighbouring F:G columns.
$ws.Columns.Item(8).ColumnWidth = $refWidth
$ws.Columns.Item(9).ColumnWidth = $refWidth

# --- Row 10 (data placeholders) : fill the new H10 / I10 cells ---
# Setting these before the row-8 header labels so the shared-string table
# picks up the same ordering as the target workbook.
$ws.Cells.Item(10, 8).Value = "`${vo.ticketPrice}"
$ws.Cells.Item(10, 9).Value = "`${vo.revenue}"

# --- Row 8 (header) : fill the new H8 / I8 header labels ---
$ws.Cells.Item(8, 8).Value = "Số tiền trước CK"
$ws.Cells.Item(8, 9).Value = "Doanh thu sau CK"

# Normalise G10/H10/I10 formatting: copy the format already used by the
# neighbouring data cell (C10) so they all share a single style record
# instead of leaving behind the old, now-redundant one.
$ws.Cells.Item(10, 3).Copy()
$ws.Cells.Item(10, 7).PasteSpecial(-4122)
$ws.Cells.Item(10, 8).PasteSpecial(-4122)
$ws.Cells.Item(10, 9).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 12 (totals bar) : the "${total}" placeholder visually moves from
# G12 (now just a blank banner cell) to the new I12 cell. ---
$totalPlaceholder = $ws.Cells.Item(12, 7).Value2
$ws.Cells.Item(12, 7).ClearContents()
$ws.Cells.Item(12, 9).Value = $totalPlaceholder

# Update the selection to match the post-edit cursor position.
$ws.Cells.Item(10, 10).Select()
